$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new column BB values first (plain numbers, no formatting side effects)
$ws.Range("BB1").Value = 45986
$ws.Range("BB3").Value = 1.457481710491582
$ws.Range("BB4").Value = 1.816525349942233
$ws.Range("BB5").Value = 5.154576318286241
$ws.Range("BB6").Value = 3.237295655668815
$ws.Range("BB7").Value = 0.9877679306576237
$ws.Range("BB8").Value = 2.422799116969476
$ws.Range("BB9").Value = 2.851134222122798
$ws.Range("BB10").Value = 2.51345851603284
$ws.Range("BB11").Value = 3.529772839278777
$ws.Range("BB12").Value = 2.406741336461815
$ws.Range("BB13").Value = 3.06399603345644
$ws.Range("BB14").Value = 0.6221001279347327
$ws.Range("BB15").Value = -1.005670135925762
$ws.Range("BB16").Value = 2.512564721370358
$ws.Range("BB17").Value = 0.6556783400334387
$ws.Range("BB18").Value = 1.751198481007687
$ws.Range("BB19").Value = 2.46481303148316
$ws.Range("BB20").Value = 2.928470412166684
$ws.Range("BB21").Value = 2.439164702314445

# Copy the header cell's formatting (date style) from BA1 onto BB1 without
# disturbing the value we just wrote.
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null
